# Commit: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to describe a generic "Property" table is being
# renamed to "DataNode" (matching the new shared naming scheme used by
# DataNode / DataTable / Entity configs), and the author's last selected
# cell before saving was C36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the only worksheet from "Property1" to "DataNode".
$ws.Name = "DataNode"

# Restore/record the author's selection at save time.
$ws.Range("C36").Select()
